$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Anxa1"
$ws.Cells.Item(2, 3).Value = "Fpr2"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 21.39646633333334
$ws.Cells.Item(2, 8).Value = 64.18939900000001
$ws.Cells.Item(2, 9).Value = 0.0721325008796955
$ws.Cells.Item(2, 10).Value = 0.0721325008796955
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.257727
$ws.Cells.Item(2, 14).Value = 3.773181
$ws.Cells.Item(2, 15).Value = 0.3633293041311343
$ws.Cells.Item(2, 16).Value = 0.3633293041311343
$ws.Cells.Item(2, 17).Value = 26.91091341202434
$ws.Cells.Item(2, 18).Value = 242.1982207082191
$ws.Cells.Item(2, 19).Value = 0.0262078513498582
$ws.Cells.Item(2, 20).Value = 0.0262078513498582

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Anxa1"
$ws.Cells.Item(3, 3).Value = "Fpr2"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 21.39646633333334
$ws.Cells.Item(3, 8).Value = 64.18939900000001
$ws.Cells.Item(3, 9).Value = 0.0721325008796955
$ws.Cells.Item(3, 10).Value = 0.0721325008796955
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.112352666666667
$ws.Cells.Item(3, 14).Value = 6.337058
$ws.Cells.Item(3, 15).Value = 0.6102116154455982
$ws.Cells.Item(3, 16).Value = 0.6102116154455982
$ws.Cells.Item(3, 17).Value = 45.19688271646023
$ws.Cells.Item(3, 18).Value = 406.7719444481421
$ws.Cells.Item(3, 19).Value = 0.04401608988793002
$ws.Cells.Item(3, 20).Value = 0.04401608988793002

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Anxa1"
$ws.Cells.Item(4, 3).Value = "Fpr2"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 21.39646633333334
$ws.Cells.Item(4, 8).Value = 64.18939900000001
$ws.Cells.Item(4, 9).Value = 0.0721325008796955
$ws.Cells.Item(4, 10).Value = 0.0721325008796955
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.09159266666666667
$ws.Cells.Item(4, 14).Value = 0.274778
$ws.Cells.Item(4, 15).Value = 0.02645908042326748
$ws.Cells.Item(4, 16).Value = 0.02645908042326749
$ws.Cells.Item(4, 17).Value = 1.959759408713556
$ws.Cells.Item(4, 18).Value = 17.637834678422
$ws.Cells.Item(4, 19).Value = 0.001908559641907276
$ws.Cells.Item(4, 20).Value = 0.001908559641907276

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Anxa1"
$ws.Cells.Item(5, 3).Value = "Fpr2"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 130.868154
$ws.Cells.Item(5, 8).Value = 392.604462
$ws.Cells.Item(5, 9).Value = 0.441187207572817
$ws.Cells.Item(5, 10).Value = 0.441187207572817
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.257727
$ws.Cells.Item(5, 14).Value = 3.773181
$ws.Cells.Item(5, 15).Value = 0.3633293041311343
$ws.Cells.Item(5, 16).Value = 0.3633293041311343
$ws.Cells.Item(5, 17).Value = 164.596410725958
$ws.Cells.Item(5, 18).Value = 1481.367696533622
$ws.Cells.Item(5, 19).Value = 0.1602962411189899
$ws.Cells.Item(5, 20).Value = 0.1602962411189899

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Anxa1"
$ws.Cells.Item(6, 3).Value = "Fpr2"
$ws.Cells.Item(6, 4).Value = "M2"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 130.868154
$ws.Cells.Item(6, 8).Value = 392.604462
$ws.Cells.Item(6, 9).Value = 0.441187207572817
$ws.Cells.Item(6, 10).Value = 0.441187207572817
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.112352666666667
$ws.Cells.Item(6, 14).Value = 6.337058
$ws.Cells.Item(6, 15).Value = 0.6102116154455982
$ws.Cells.Item(6, 16).Value = 0.6102116154455982
$ws.Cells.Item(6, 17).Value = 276.439694083644
$ws.Cells.Item(6, 18).Value = 2487.957246752796
$ws.Cells.Item(6, 19).Value = 0.2692175586469411
$ws.Cells.Item(6, 20).Value = 0.2692175586469411

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Anxa1"
$ws.Cells.Item(7, 3).Value = "Fpr2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 130.868154
$ws.Cells.Item(7, 8).Value = 392.604462
$ws.Cells.Item(7, 9).Value = 0.441187207572817
$ws.Cells.Item(7, 10).Value = 0.441187207572817
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.09159266666666667
$ws.Cells.Item(7, 14).Value = 0.274778
$ws.Cells.Item(7, 15).Value = 0.02645908042326748
$ws.Cells.Item(7, 16).Value = 0.02645908042326749
$ws.Cells.Item(7, 17).Value = 11.986563206604
$ws.Cells.Item(7, 18).Value = 107.879068859436
$ws.Cells.Item(7, 19).Value = 0.01167340780688597
$ws.Cells.Item(7, 20).Value = 0.01167340780688597

$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Anxa1"
$ws.Cells.Item(8, 3).Value = "Fpr2"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 75.02619166666666
$ws.Cells.Item(8, 8).Value = 225.078575
$ws.Cells.Item(8, 9).Value = 0.2529308696158396
$ws.Cells.Item(8, 10).Value = 0.2529308696158397
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 1.257727
$ws.Cells.Item(8, 14).Value = 3.773181
$ws.Cells.Item(8, 15).Value = 0.3633293041311343
$ws.Cells.Item(8, 16).Value = 0.3633293041311343
$ws.Cells.Item(8, 17).Value = 94.36246696634166
$ws.Cells.Item(8, 18).Value = 849.262202697075
$ws.Cells.Item(8, 19).Value = 0.09189719685080566
$ws.Cells.Item(8, 20).Value = 0.09189719685080569

$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Anxa1"
$ws.Cells.Item(9, 3).Value = "Fpr2"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 75.02619166666666
$ws.Cells.Item(9, 8).Value = 225.078575
$ws.Cells.Item(9, 9).Value = 0.2529308696158396
$ws.Cells.Item(9, 10).Value = 0.2529308696158397
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.112352666666667
$ws.Cells.Item(9, 14).Value = 6.337058
$ws.Cells.Item(9, 15).Value = 0.6102116154455982
$ws.Cells.Item(9, 16).Value = 0.6102116154455982
$ws.Cells.Item(9, 17).Value = 158.4817760369278
$ws.Cells.Item(9, 18).Value = 1426.33598433235
$ws.Cells.Item(9, 19).Value = 0.1543413545443414
$ws.Cells.Item(9, 20).Value = 0.1543413545443415

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Anxa1"
$ws.Cells.Item(10, 3).Value = "Fpr2"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 75.02619166666666
$ws.Cells.Item(10, 8).Value = 225.078575
$ws.Cells.Item(10, 9).Value = 0.2529308696158396
$ws.Cells.Item(10, 10).Value = 0.2529308696158397
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.09159266666666667
$ws.Cells.Item(10, 14).Value = 0.274778
$ws.Cells.Item(10, 15).Value = 0.02645908042326748
$ws.Cells.Item(10, 16).Value = 0.02645908042326749
$ws.Cells.Item(10, 17).Value = 6.871848964594444
$ws.Cells.Item(10, 18).Value = 61.84664068135
$ws.Cells.Item(10, 19).Value = 0.006692318220692482
$ws.Cells.Item(10, 20).Value = 0.006692318220692484

$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Anxa1"
$ws.Cells.Item(11, 3).Value = "Fpr2"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 69.33645133333333
$ws.Cells.Item(11, 8).Value = 208.009354
$ws.Cells.Item(11, 9).Value = 0.2337494219316478
$ws.Cells.Item(11, 10).Value = 0.2337494219316478
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.257727
$ws.Cells.Item(11, 14).Value = 3.773181
$ws.Cells.Item(11, 15).Value = 0.3633293041311343
$ws.Cells.Item(11, 16).Value = 0.3633293041311343
$ws.Cells.Item(11, 17).Value = 87.20632692611933
$ws.Cells.Item(11, 18).Value = 784.856942335074
$ws.Cells.Item(11, 19).Value = 0.08492801481148049
$ws.Cells.Item(11, 20).Value = 0.0849280148114805

$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Anxa1"
$ws.Cells.Item(12, 3).Value = "Fpr2"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 69.33645133333333
$ws.Cells.Item(12, 8).Value = 208.009354
$ws.Cells.Item(12, 9).Value = 0.2337494219316478
$ws.Cells.Item(12, 10).Value = 0.2337494219316478
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.112352666666667
$ws.Cells.Item(12, 14).Value = 6.337058
$ws.Cells.Item(12, 15).Value = 0.6102116154455982
$ws.Cells.Item(12, 16).Value = 0.6102116154455982
$ws.Cells.Item(12, 17).Value = 146.4630378711702
$ws.Cells.Item(12, 18).Value = 1318.167340840532
$ws.Cells.Item(12, 19).Value = 0.1426366123663855
$ws.Cells.Item(12, 20).Value = 0.1426366123663856

$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Anxa1"
$ws.Cells.Item(13, 3).Value = "Fpr2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 69.33645133333333
$ws.Cells.Item(13, 8).Value = 208.009354
$ws.Cells.Item(13, 9).Value = 0.2337494219316478
$ws.Cells.Item(13, 10).Value = 0.2337494219316478
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.09159266666666667
$ws.Cells.Item(13, 14).Value = 0.274778
$ws.Cells.Item(13, 15).Value = 0.02645908042326748
$ws.Cells.Item(13, 16).Value = 0.02645908042326749
$ws.Cells.Item(13, 17).Value = 6.350710474823555
$ws.Cells.Item(13, 18).Value = 57.156394273412
$ws.Cells.Item(13, 19).Value = 0.006184794753781753
$ws.Cells.Item(13, 20).Value = 0.006184794753781754
